# Task panel bugfixes: enable task review / save / load.
# Insert 4 new eventList rows right after the existing "taskInfo" row (row 215)
# and wire "taskInfo" to the new "taskInfoHasTask" branch instead of the old
# no-op (";") value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the 4 new rows (216:219); existing rows 216+ shift down to 220+.
$ws.Rows("216:219").Insert()

# taskInfo now routes through closeWindow -> taskInfoHasTask instead of doing nothing.
$ws.Range("D215").Value = "closeWindow;taskInfoHasTask"

# New row 216: taskInfoHasTask condition - branches on whether a task is currently active.
$ws.Range("A216").Value = "taskInfoHasTask"
$ws.Range("B216").Value = "接取了任务"
$ws.Range("C216").Value = "condition"
$ws.Range("D216").Value = "hasTask;taskInfoShowCurrentTask;taskInfoNoTask"

# New row 217: taskInfoShowCurrentTask window - shows the task info panel.
$ws.Range("A217").Value = "taskInfoShowCurrentTask"
$ws.Range("C217").Value = "window"
$ws.Range("D217").Value = "TaskPanel;infoList;infoList;1"

# New row 218: taskInfoNoTask eventList - no active task, show the "no task" dialog.
$ws.Range("A218").Value = "taskInfoNoTask"
$ws.Range("C218").Value = "eventList"
$ws.Range("D218").Value = "taskInfoDialogNoTask;infoList"

# New row 219: taskInfoDialogNoTask dialog - the actual "no task received" message.
$ws.Range("A219").Value = "taskInfoDialogNoTask"
$ws.Range("C219").Value = "dialog"
$ws.Range("D219").Value = "dialog_no_task_received"
